# Updated cryptos list (Price / Volume(1h) refresh, plus a couple of
# ranking swaps: Uniswap<->ImmutableX at rows 19/20, RenderToken<->Celestia
# at rows 38/39). All D/E cells are text in this sheet, so numeric-looking
# values are written with NumberFormat "@" (then reset to the default
# "Normal" style) to stop Excel's COM layer from auto-converting them to
# real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.967.78"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "2.822.13"
$ws.Range("E3").Value = "  +3.04%  "
$ws.Range("E4").Value = "  +0.12%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "354.91"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +6.61%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "113.42"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.94%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.546"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.604"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +6.57%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "41.93"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.58%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0845"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.55%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "20.05"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.05%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.130"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.32%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.75"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").Value = "3.253.33"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").Value = "2.830.68"
$ws.Range("E16").Value = "  +2.21%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.893"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.43%  "
$ws.Range("D18").Value = "52.023.79"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.16"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.25"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +6.60%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.74"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("D22").Value = "0.0₃0992"
$ws.Range("E22").Value = "  +2.33%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "269.67"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -3.30%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "69.62"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("E25").Value = "  +5.90%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "26.66"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  +0.12%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.30"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("E29").Value = "  +1.37%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.140"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.66%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "50.71"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.63%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "33.78"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.67%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.86"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +6.07%  "
$ws.Range("E34").Value = "  +28.87%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0830"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "18.40"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.72%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "4.86"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.62%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.20"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.07%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.59"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +8.62%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "127.71"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.24%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "23.40"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.15%  "
$ws.Range("E44").Value = "  +1.86%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.30"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.46%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.34"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").Value = "2.049.83"
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("E48").Value = "  +3.83%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.953"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +10.70%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "5.68"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.08%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "60.25"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.26%  "
